$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: update title and link
$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"
$ws.Range("E26").Value = "https://blog.est.ai/2022/06/%ec%9d%b8%ea%b3%b5%ec%a7%80%eb%8a%a5-%ec%9d%8c%ec%84%b1-%ec%83%9d%ec%84%b1-%ec%97%b0%ea%b5%ac-%ec%9d%8c%ec%84%b1-%eb%b6%84%eb%a5%98-%ec%86%94%eb%a3%a8%ec%85%98/"

# Row 36: update title and link
$ws.Range("D36").Value = "Toward Optimal Optimizer"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/368"

# Row 46: update title and link
$ws.Range("D46").Value = "요로감염증 (UTI)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/478"
